$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric values per the diff
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5

# Move the active selection from C4 to F4
$ws.Range("F4").Select()
